$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tugas")

# ---- Row 1 / Row 2 headers ----
$ws.Range("A1").Value = "Data"
$ws.Range("B1").Value = "Waktu Sebelum Tunning (ms)"
$ws.Range("G1").Value = "Waktu Sesudah Tuning (ms)"

$ws.Range("B2").Value = "QUERY 1"
$ws.Range("C2").Value = "QUERY 2"
$ws.Range("D2").Value = "QUERY 3"
$ws.Range("E2").Value = "QUERY 4"
$ws.Range("F2").Value = "QUERY 5"
$ws.Range("G2").Value = "QUERY 1"
$ws.Range("H2").Value = "QUERY 2"
$ws.Range("I2").Value = "QUERY 3"
$ws.Range("J2").Value = "QUERY 4"
$ws.Range("K2").Value = "QUERY 5"

# ---- Data rows 3-6 (numeric timings) ----
$ws.Range("B3").Value = 0.057710299999999999
$ws.Range("C3").Value = 0.0003678
$ws.Range("D3").Value = 0.00034909999999999997
$ws.Range("E3").Value = 0.091142500000000001
$ws.Range("F3").Value = 0.050515299999999999
$ws.Range("G3").Value = 0.00041110000000000002
$ws.Range("H3").Value = 0.00034600000000000001
$ws.Range("I3").Value = 0.00034309999999999999
$ws.Range("J3").Value = 0.00016470999999999999
$ws.Range("K3").Value = 0.00014311

$ws.Range("B4").Value = 0.00046359999999999999
$ws.Range("C4").Value = 0.00062620000000000004
$ws.Range("D4").Value = 0.0036143
$ws.Range("E4").Value = 0.0030785000000000001
$ws.Range("F4").Value = 0.00052360000000000004
$ws.Range("G4").Value = 0.0004013
$ws.Range("H4").Value = 0.00047390000000000003
$ws.Range("I4").Value = 0.00040840000000000001
$ws.Range("J4").Value = 0.0037843
$ws.Range("K4").Value = 0.0029432

$ws.Range("B5").Value = 0.082456600000000005
$ws.Range("C5").Value = 0.00092960000000000004
$ws.Range("D5").Value = 0.00050060000000000002
$ws.Range("E5").Value = 0.53413699999999997
$ws.Range("F5").Value = 0.56575640000000005
$ws.Range("G5").Value = 0.00061359999999999995
$ws.Range("H5").Value = 0.1314535
$ws.Range("I5").Value = 0.0005042
$ws.Range("J5").Value = 0.62759430000000005
$ws.Range("K5").Value = 0.33083810000000002

$ws.Range("B6").Value = 0.044920399999999999
$ws.Range("C6").Value = 0.00033159999999999998
$ws.Range("D6").Value = 0.00036388000000000002
$ws.Range("E6").Value = 0.0015384999999999999
$ws.Range("F6").Value = 0.0026548000000000001
$ws.Range("G6").Value = 0.00034519999999999999
$ws.Range("H6").Value = 0.00034870000000000002
$ws.Range("I6").Value = 0.00031629999999999999
$ws.Range("J6").Value = 0.0097050999999999995
$ws.Range("K6").Value = 0.0014040000000000001

# ---- Rows 7-9 dash placeholders ----
$ws.Range("B7:K7").Value = "-"
$ws.Range("B8:K8").Value = "-"
$ws.Range("B9:K9").Value = "-"

# ---- Borders: full thin grid border for the whole A1:K9 data table ----
$rng = $ws.Range("A1:K9")
$rng.Borders.LineStyle = 1
$rng.Borders.Item(11).LineStyle = 1
$rng.Borders.Item(12).LineStyle = 1

# ---- Bold header rows (1-2) ----
$ws.Range("A1:K2").Font.Bold = $true

# ---- Alignment ----
$ws.Range("A1:K1").HorizontalAlignment = -4108
$ws.Range("B1:K1").VerticalAlignment = -4108
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4108
$ws.Range("B7:K9").HorizontalAlignment = -4108

# ---- Merges ----
$ws.Range("B1:F1").Merge()
$ws.Range("G1:K1").Merge()
$ws.Range("A1:A2").Merge()

# ---- Column widths ----
$ws.Range("B1:F1").ColumnWidth = 8.14
$ws.Range("G1").ColumnWidth = 8.4
$ws.Range("H1").ColumnWidth = 7.76
$ws.Range("I1").ColumnWidth = 7.63
$ws.Range("J1:K1").ColumnWidth = 8.14
